$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all values in A1:H14 from 0.5 to 3.048E-5
$ws.Range("A1:H14").Value = 0.00003048

# Column A now needs a best-fit width to show the new shorter values (width=11 in the saved file).
# The host engine's ColumnWidth setter pads by ~0.83 chars internally, so back that out here
# so the persisted <col> width lands on exactly 11.
$ws.Columns.Item(1).ColumnWidth = 10.1667
